$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for "Uridine" at row 21 (pushes existing rows 21-27 down to 22-28)
$ws.Rows.Item(21).Insert()
$ws.Range("A21").Value = "Uridine"
$ws.Range("B21").Value = 245.0768
$ws.Range("C21").Value = 31

# Re-sort the full data range (A2:C28) by column B, ascending, matching the
# "Data > Sort" operation the author performed (produces the <sortState> the
# workbook now carries).
$srt = $ws.Sort
$srt.SortFields.Clear()
$srt.SortFields.Add($ws.Range("B2:B28"))
$srt.SetRange($ws.Range("A1:C28"))
$srt.Header = 1
$srt.Apply()

# Update the active selection to match the author's last recorded cursor position.
$ws.Range("G10").Select()
